$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Notes")

# Update the Description text (row 2)
$ws.Range("A2").Value = "Description: Household Sanitation Coverage (%)"

# Update the Source text (row 4)
$ws.Range("A4").Value = "Source: Water and Environment Sector Performance Reports 2010-2014 - Ministry of Water and Environment"

# Insert a new row after the Source row for the Source-link
$ws.Rows("5").Insert()
$ws.Range("A5").Value = "Source-link: http://www.mwe.go.ug/index.php?option=com_docman&task=cat_view&Itemid=223&gid=15"

# Update the license text (now on row 14, after the inserted row shifted things down)
$ws.Range("A14").Value = "It is licensed under a Creative Commons Attribution 4.0 International license."

# Insert a new row after the license row with a link to more licensing info
$ws.Rows("15").Insert()
$ws.Range("A15").Value = "More information on licensing is available here: https://creativecommons.org/licenses/by/4.0/"
